# Add new RATE cards and update ID for Delta the Magnet Warrior
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RATE-JP")

# Fill in the ygopro ids (column B) for all card rows (2-81).
# Row 2 starts at 100911001 and increments by 1 per row.
for ($r = 2; $r -le 81; $r++) {
    $ws.Cells.Item($r, 2).Value = 100911000 + ($r - 1)
}

# New cards added to the list (column A names, referencing new shared strings).
$ws.Range("A7").Value = "Speedroid Bamboo Horse"
$ws.Range("A8").Value = "Wind Witch - Ice Bell"
$ws.Range("A70").Value = "Cipher Spectrum"

# The sheet previously had one extra trailing blank row (row 82); remove it
# so the data now ends at row 81.
$ws.Rows.Item(82).Delete()

# Update the selected/active cell in the sheet view.
$ws.Range("B9").Select()
